$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "p1 <= p2,p3 <= p2,p3 <= p4"
$ws.Range("B5").Value = "[0 0 0 0]"
$ws.Range("C5").Value = "T>R>P>S"
$ws.Range("D5").Value = "Behave better on bad conditions, no reason to cooperate."

$ws.Range("A5:D5").HorizontalAlignment = -4108
$ws.Range("A5:D5").VerticalAlignment = -4108

$ws.Range("A5").Select()
